$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.613107666666667
$ws.Range("H2").Value = 22.839323
$ws.Range("I2").Value = 0.08102996839946881
$ws.Range("J2").Value = 0.0810299683994688
$ws.Range("M2").Value = 0.5373756666666667
$ws.Range("N2").Value = 1.612127
$ws.Range("O2").Value = 0.007472820128982582
$ws.Range("P2").Value = 0.007472820128982581
$ws.Range("Q2").Value = 4.091098807780112
$ws.Range("R2").Value = 36.819889270021
$ws.Range("S2").Value = 0.0006055223789063731
$ws.Range("T2").Value = 0.0006055223789063729
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.613107666666667
$ws.Range("H3").Value = 22.839323
$ws.Range("I3").Value = 0.08102996839946881
$ws.Range("J3").Value = 0.0810299683994688
$ws.Range("O3").Value = 0.1537223653287423
$ws.Range("P3").Value = 0.1537223653287423
$ws.Range("Q3").Value = 84.15743650599279
$ws.Range("R3").Value = 757.4169285539351
$ws.Range("S3").Value = 0.01245611840487959
$ws.Range("T3").Value = 0.01245611840487959
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.613107666666667
$ws.Range("H4").Value = 22.839323
$ws.Range("I4").Value = 0.08102996839946881
$ws.Range("J4").Value = 0.0810299683994688
$ws.Range("M4").Value = 30.561198
$ws.Range("N4").Value = 91.683594
$ws.Range("O4").Value = 0.4249882340167162
$ws.Range("P4").Value = 0.4249882340167161
$ws.Range("Q4").Value = 232.665690796318
$ws.Range("R4").Value = 2093.991217166862
$ws.Range("S4").Value = 0.03443678317252057
$ws.Range("T4").Value = 0.03443678317252056
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.613107666666667
$ws.Range("H5").Value = 22.839323
$ws.Range("I5").Value = 0.08102996839946881
$ws.Range("J5").Value = 0.0810299683994688
$ws.Range("M5").Value = 29.75783666666667
$ws.Range("N5").Value = 89.27351
$ws.Range("O5").Value = 0.4138165805255589
$ws.Range("P5").Value = 0.4138165805255589
$ws.Range("Q5").Value = 226.5496144704144
$ws.Range("R5").Value = 2038.94653023373
$ws.Range("S5").Value = 0.03353154444316228
$ws.Range("T5").Value = 0.03353154444316227
$ws.Range("I6").Value = 0.7831116101658118
$ws.Range("J6").Value = 0.7831116101658117
$ws.Range("M6").Value = 0.5373756666666667
$ws.Range("N6").Value = 1.612127
$ws.Range("O6").Value = 0.007472820128982582
$ws.Range("P6").Value = 0.007472820128982581
$ws.Range("Q6").Value = 39.5382972249699
$ws.Range("R6").Value = 355.8446750247291
$ws.Range("S6").Value = 0.005852052203687039
$ws.Range("T6").Value = 0.005852052203687038
$ws.Range("I7").Value = 0.7831116101658118
$ws.Range("J7").Value = 0.7831116101658117
$ws.Range("O7").Value = 0.1537223653287423
$ws.Range("P7").Value = 0.1537223653287423
$ws.Range("S7").Value = 0.1203817690310886
$ws.Range("T7").Value = 0.1203817690310885
$ws.Range("I8").Value = 0.7831116101658118
$ws.Range("J8").Value = 0.7831116101658117
$ws.Range("M8").Value = 30.561198
$ws.Range("N8").Value = 91.683594
$ws.Range("O8").Value = 0.4249882340167162
$ws.Range("P8").Value = 0.4249882340167161
$ws.Range("Q8").Value = 2248.590334524182
$ws.Range("R8").Value = 20237.31301071764
$ws.Range("S8").Value = 0.3328132202423554
$ws.Range("T8").Value = 0.3328132202423554
$ws.Range("I9").Value = 0.7831116101658118
$ws.Range("J9").Value = 0.7831116101658117
$ws.Range("M9").Value = 29.75783666666667
$ws.Range("N9").Value = 89.27351
$ws.Range("O9").Value = 0.4138165805255589
$ws.Range("P9").Value = 0.4138165805255589
$ws.Range("Q9").Value = 2189.481705037086
$ws.Range("R9").Value = 19705.33534533377
$ws.Range("S9").Value = 0.3240645686886808
$ws.Range("T9").Value = 0.3240645686886807
$ws.Range("G10").Value = 12.72068066666667
$ws.Range("H10").Value = 38.162042
$ws.Range("I10").Value = 0.1353923256534006
$ws.Range("J10").Value = 0.1353923256534005
$ws.Range("M10").Value = 0.5373756666666667
$ws.Range("N10").Value = 1.612127
$ws.Range("O10").Value = 0.007472820128982582
$ws.Range("P10").Value = 0.007472820128982581
$ws.Range("Q10").Value = 6.835784253703778
$ws.Range("R10").Value = 61.522058283334
$ws.Range("S10").Value = 0.001011762496452497
$ws.Range("T10").Value = 0.001011762496452496
$ws.Range("G11").Value = 12.72068066666667
$ws.Range("H11").Value = 38.162042
$ws.Range("I11").Value = 0.1353923256534006
$ws.Range("J11").Value = 0.1353923256534005
$ws.Range("O11").Value = 0.1537223653287423
$ws.Range("P11").Value = 0.1537223653287423
$ws.Range("Q11").Value = 140.6179870810545
$ws.Range("R11").Value = 1265.56188372949
$ws.Range("S11").Value = 0.02081282854680009
$ws.Range("T11").Value = 0.02081282854680009
$ws.Range("G12").Value = 12.72068066666667
$ws.Range("H12").Value = 38.162042
$ws.Range("I12").Value = 0.1353923256534006
$ws.Range("J12").Value = 0.1353923256534005
$ws.Range("M12").Value = 30.561198
$ws.Range("N12").Value = 91.683594
$ws.Range("O12").Value = 0.4249882340167162
$ws.Range("P12").Value = 0.4249882340167161
$ws.Range("Q12").Value = 388.759240548772
$ws.Range("R12").Value = 3498.833164938948
$ws.Range("S12").Value = 0.05754014537885484
$ws.Range("T12").Value = 0.05754014537885483
$ws.Range("G13").Value = 12.72068066666667
$ws.Range("H13").Value = 38.162042
$ws.Range("I13").Value = 0.1353923256534006
$ws.Range("J13").Value = 0.1353923256534005
$ws.Range("M13").Value = 29.75783666666667
$ws.Range("N13").Value = 89.27351
$ws.Range("O13").Value = 0.4138165805255589
$ws.Range("P13").Value = 0.4138165805255589
$ws.Range("Q13").Value = 378.5399375674911
$ws.Range("R13").Value = 3406.85943810742
$ws.Range("S13").Value = 0.05602758923129313
$ws.Range("T13").Value = 0.05602758923129312
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.04379166666666667
$ws.Range("H14").Value = 0.131375
$ws.Range("I14").Value = 0.0004660957813189216
$ws.Range("J14").Value = 0.0004660957813189215
$ws.Range("M14").Value = 0.5373756666666667
$ws.Range("N14").Value = 1.612127
$ws.Range("O14").Value = 0.007472820128982582
$ws.Range("P14").Value = 0.007472820128982581
$ws.Range("Q14").Value = 0.02353257606944445
$ws.Range("R14").Value = 0.211793184625
$ws.Range("S14").Value = 0.000003483049936673901
$ws.Range("T14").Value = 0.0000034830499366739
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.04379166666666667
$ws.Range("H15").Value = 0.131375
$ws.Range("I15").Value = 0.0004660957813189216
$ws.Range("J15").Value = 0.0004660957813189215
$ws.Range("O15").Value = 0.1537223653287423
$ws.Range("P15").Value = 0.1537223653287423
$ws.Range("Q15").Value = 0.4840854179861113
$ws.Range("R15").Value = 4.356768761875001
$ws.Range("S15").Value = 0.00007164934597409286
$ws.Range("T15").Value = 0.00007164934597409284
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.04379166666666667
$ws.Range("H16").Value = 0.131375
$ws.Range("I16").Value = 0.0004660957813189216
$ws.Range("J16").Value = 0.0004660957813189215
$ws.Range("M16").Value = 30.561198
$ws.Range("N16").Value = 91.683594
$ws.Range("O16").Value = 0.4249882340167162
$ws.Range("P16").Value = 0.4249882340167161
$ws.Range("Q16").Value = 1.33832579575
$ws.Range("R16").Value = 12.04493216175
$ws.Range("S16").Value = 0.00019808522298537
$ws.Range("T16").Value = 0.00019808522298537
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.04379166666666667
$ws.Range("H17").Value = 0.131375
$ws.Range("I17").Value = 0.0004660957813189216
$ws.Range("J17").Value = 0.0004660957813189215
$ws.Range("M17").Value = 29.75783666666667
$ws.Range("N17").Value = 89.27351
$ws.Range("O17").Value = 0.4138165805255589
$ws.Range("P17").Value = 0.4138165805255589
$ws.Range("Q17").Value = 1.303145264027778
$ws.Range("R17").Value = 11.72830737625
$ws.Range("S17").Value = 0.0001928781624227848
$ws.Range("T17").Value = 0.0001928781624227848
